$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 345.77
$ws.Range("I15").Value = 345.77
$ws.Range("K15").Value = 1037.31
$ws.Range("M15").Value = -868.3099999999999
$ws.Range("H28").Value = 691.9545000000001
$ws.Range("J28").Value = 1201.25
$ws.Range("L28").Value = 1201.25
$ws.Range("N28").Value = -2171.25
$ws.Range("H74").Value = 2977.875
$ws.Range("I74").Value = 2720.7693
$ws.Range("J74").Value = 3281.7273
$ws.Range("K74").Value = 2720.7693
$ws.Range("L74").Value = 3281.7273
$ws.Range("M74").Value = -1784.7693
$ws.Range("N74").Value = -5153.7273
$ws.Range("H77").Value = 2977.875
$ws.Range("I77").Value = 2720.7693
$ws.Range("J77").Value = 3281.7273
$ws.Range("K77").Value = 13603.8465
$ws.Range("L77").Value = 16408.6365
$ws.Range("M77").Value = -8923.8465
$ws.Range("N77").Value = -25768.6365
$ws.Range("H129").Value = 2807.08
$ws.Range("I129").Value = 8226.691999999999
$ws.Range("J129").Value = 902.8919
$ws.Range("K129").Value = 24680.076
$ws.Range("L129").Value = 2708.6757
$ws.Range("M129").Value = -19680.076
$ws.Range("N129").Value = -12708.6757
$ws.Range("H135").Value = 1501.3948
$ws.Range("I135").Value = 464.36365
$ws.Range("J135").Value = 2927.3125
$ws.Range("K135").Value = 4179.27285
$ws.Range("L135").Value = 26345.8125
$ws.Range("M135").Value = -1644.27285
$ws.Range("N135").Value = -31415.8125
$ws.Range("H138").Value = 2549.2693
$ws.Range("I138").Value = 1766.25
$ws.Range("J138").Value = 2897.2778
$ws.Range("K138").Value = 5298.75
$ws.Range("L138").Value = 8691.8334
$ws.Range("M138").Value = -158.75
$ws.Range("N138").Value = -18971.8334
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27317.36
$ws.Range("I32").Value = 6668.788
$ws.Range("K32").Value = 6668.788
$ws.Range("M32").Value = -6381.788
$ws.Range("H61").Value = 1497.2295
$ws.Range("I61").Value = 1032.2
$ws.Range("J61").Value = 2805.125
$ws.Range("K61").Value = 1032.2
$ws.Range("L61").Value = 2805.125
$ws.Range("M61").Value = -820.2
$ws.Range("N61").Value = -3229.125
$ws.Range("H63").Value = 2214.1667
$ws.Range("I63").Value = 2084.889
$ws.Range("J63").Value = 2602
$ws.Range("K63").Value = 2084.889
$ws.Range("L63").Value = 2602
$ws.Range("M63").Value = -1398.889
$ws.Range("N63").Value = -3974
$ws.Range("H66").Value = 2214.1667
$ws.Range("I66").Value = 2084.889
$ws.Range("J66").Value = 2602
$ws.Range("K66").Value = 10424.445
$ws.Range("L66").Value = 13010
$ws.Range("M66").Value = -6992.445
$ws.Range("N66").Value = -19874
$ws.Range("H74").Value = 836
$ws.Range("I74").Value = 809.38464
$ws.Range("K74").Value = 809.38464
$ws.Range("M74").Value = 64.61536000000001
$ws.Range("H77").Value = 836
$ws.Range("I77").Value = 809.38464
$ws.Range("K77").Value = 4046.9232
$ws.Range("M77").Value = 321.0767999999998
$ws.Range("H97").Value = 24892.262
$ws.Range("I97").Value = 34006.1
$ws.Range("K97").Value = 34006.1
$ws.Range("M97").Value = -33510.1
$ws.Range("H132").Value = 12966.84
$ws.Range("I132").Value = 16552.75
$ws.Range("K132").Value = 49658.25
$ws.Range("M132").Value = -47128.25
$ws.Range("H136").Value = 1497.2295
$ws.Range("I136").Value = 1032.2
$ws.Range("J136").Value = 2805.125
$ws.Range("K136").Value = 3096.6
$ws.Range("L136").Value = 8415.375
$ws.Range("M136").Value = -546.6000000000004
$ws.Range("N136").Value = -13515.375
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2117.0881
$ws.Range("I134").Value = 1395.7667
$ws.Range("J134").Value = 7527
$ws.Range("K134").Value = 4187.300099999999
$ws.Range("L134").Value = 22581
$ws.Range("M134").Value = -1652.300099999999
$ws.Range("N134").Value = -27651
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4721.353
$ws.Range("I132").Value = 5033.727
$ws.Range("J132").Value = 4148.6665
$ws.Range("K132").Value = 15101.181
$ws.Range("L132").Value = 12445.9995
$ws.Range("M132").Value = -12571.181
$ws.Range("N132").Value = -17505.9995
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1222.2222
$ws.Range("J34").Value = 1350
$ws.Range("L34").Value = 4050
$ws.Range("N34").Value = -4218
$ws.Range("H58").Value = 2400
$ws.Range("J58").Value = 2400
$ws.Range("L58").Value = 7200
$ws.Range("N58").Value = -7456
$ws.Range("H113").Value = 1009.5238
$ws.Range("I113").Value = 1736.25
$ws.Range("J113").Value = 562.3077
$ws.Range("K113").Value = 5208.75
$ws.Range("L113").Value = 1686.9231
$ws.Range("M113").Value = -3038.75
$ws.Range("N113").Value = -6026.9231
$ws.Range("H120").Value = 374437.66
$ws.Range("I120").Value = 374437.66
$ws.Range("K120").Value = 1123312.98
$ws.Range("M120").Value = -1118474.98
$ws.Range("H131").Value = 1237.5555
$ws.Range("J131").Value = 1273.7922
$ws.Range("L131").Value = 3821.376600000001
$ws.Range("N131").Value = -13901.3766
$ws.Range("H137").Value = 3886592.2
$ws.Range("I137").Value = 72824.21000000001
$ws.Range("J137").Value = 8335988.5
$ws.Range("K137").Value = 218472.63
$ws.Range("L137").Value = 25007965.5
$ws.Range("M137").Value = -213372.63
$ws.Range("N137").Value = -25018165.5
$ws.Range("H141").Value = 18975
$ws.Range("I141").Value = 33120
$ws.Range("J141").Value = 4830
$ws.Range("K141").Value = 99360
$ws.Range("L141").Value = 14490
$ws.Range("M141").Value = -94180
$ws.Range("N141").Value = -24850
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1688.3226
$ws.Range("I122").Value = 1501.6666
$ws.Range("J122").Value = 1946.7693
$ws.Range("K122").Value = 4504.9998
$ws.Range("L122").Value = 5840.3079
$ws.Range("M122").Value = -2054.9998
$ws.Range("N122").Value = -10740.3079
$ws.Range("H132").Value = 2524.9773
$ws.Range("I132").Value = 1957.2
$ws.Range("K132").Value = 5871.6
$ws.Range("M132").Value = -3341.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 479814.78
$ws.Range("I16").Value = 92027.73
$ws.Range("J16").Value = 835286.25
$ws.Range("K16").Value = 92027.73
$ws.Range("L16").Value = 835286.25
$ws.Range("M16").Value = -91857.73
$ws.Range("N16").Value = -835626.25
$ws.Range("H122").Value = 2418.8
$ws.Range("I122").Value = 2381.6086
$ws.Range("J122").Value = 2846.5
$ws.Range("K122").Value = 7144.825800000001
$ws.Range("L122").Value = 8539.5
$ws.Range("M122").Value = -4694.825800000001
$ws.Range("N122").Value = -13439.5
$ws.Range("H132").Value = 4882.304
$ws.Range("I132").Value = 5172.6
$ws.Range("J132").Value = 4338
$ws.Range("K132").Value = 15517.8
$ws.Range("L132").Value = 13014
$ws.Range("M132").Value = -12987.8
$ws.Range("N132").Value = -18074
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 16400
$ws.Range("J34").Value = 16400
$ws.Range("L34").Value = 16400
$ws.Range("N34").Value = -16806
$ws.Range("H37").Value = 14016.667
$ws.Range("I37").Value = 9800
$ws.Range("J37").Value = 14860
$ws.Range("K37").Value = 9800
$ws.Range("L37").Value = 14860
$ws.Range("N37").Value = -15266
$ws.Range("M37").Value = -9597
$ws.Range("H43").Value = 16400
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 16400
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 16400
$ws.Range("N43").Value = -16698
$ws.Range("M43").ClearContents()
$ws.Range("H132").Value = 1939.8529
$ws.Range("I132").Value = 2046.94
$ws.Range("J132").Value = 1642.3889
$ws.Range("K132").Value = 6140.82
$ws.Range("L132").Value = 4927.1667
$ws.Range("M132").Value = -3610.82
$ws.Range("N132").Value = -9987.1667
$ws.Range("H136").Value = 1135.7179
$ws.Range("I136").Value = 314.8409
$ws.Range("J136").Value = 2198.0293
$ws.Range("K136").Value = 944.5227
$ws.Range("L136").Value = 6594.0879
$ws.Range("M136").Value = 1605.4773
$ws.Range("N136").Value = -11694.0879
